# Daily attendance processing - rotate "Recorded By" (column G) name lists
# so that the first listed name moves to the end of the comma-separated list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 157

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    $text = [string]$val
    if ($text -eq "") {
        continue
    }

    if ($text -notmatch ",") {
        continue
    }

    $parts = $text -split ",\s*"
    if ($parts.Count -le 1) {
        continue
    }

    $rotated = @($parts[1..($parts.Count - 1)]) + @($parts[0])
    $newText = [string]::Join(", ", $rotated)

    $cell.Value = $newText
}
